$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "302.71"
Set-TextValue "E2" "-0.18%"
Set-TextValue "G2" "15"
Set-TextValue "D3" "32.67"
Set-TextValue "E3" "1.91%"
Set-TextValue "G3" "15"
Set-TextValue "D4" "5.077"
Set-TextValue "E4" "-0.81%"
Set-TextValue "G4" "15"
Set-TextValue "D5" "0.07750"
Set-TextValue "E5" "-1.45%"
Set-TextValue "G5" "15"
Set-TextValue "D6" "2.092"
Set-TextValue "E6" "-8.60%"
Set-TextValue "G6" "15"
Set-TextValue "D7" "7.913"
Set-TextValue "E7" "1.16%"
Set-TextValue "G7" "15"
Set-TextValue "D8" "0.9269"
Set-TextValue "E8" "-0.33%"
Set-TextValue "G8" "15"
Set-TextValue "D9" "0.1776"
Set-TextValue "E9" "0.38%"
Set-TextValue "G9" "15"
Set-TextValue "D10" "0.08042"
Set-TextValue "E10" "3.02%"
Set-TextValue "G10" "15"
Set-TextValue "D11" "0.08580"
Set-TextValue "E11" "-3.42%"
Set-TextValue "G11" "15"
Set-TextValue "D12" "0.03098"
Set-TextValue "E12" "0.33%"
Set-TextValue "G12" "15"
Set-TextValue "D13" "0.09977"
Set-TextValue "E13" "-0.31%"
Set-TextValue "G13" "15"
Set-TextValue "D14" "0.001525"
Set-TextValue "E14" "0.81%"
Set-TextValue "G14" "15"
Set-TextValue "D15" "0.005655"
Set-TextValue "E15" "-3.41%"
Set-TextValue "G15" "15"
Set-TextValue "G16" "15"
Set-TextValue "E17" "0.01%"
Set-TextValue "G17" "15"
Set-TextValue "D18" "3.795"
Set-TextValue "E18" "-0.07%"
Set-TextValue "G18" "15"
Set-TextValue "D19" "2.159"
Set-TextValue "E19" "-3.98%"
Set-TextValue "G19" "15"
Set-TextValue "D20" "0.3340"
Set-TextValue "E20" "2.02%"
Set-TextValue "G20" "15"
Set-TextValue "D21" "0.1317"
Set-TextValue "E21" "-1.47%"
Set-TextValue "G21" "15"
Set-TextValue "D22" "4.399"
Set-TextValue "E22" "3.78%"
Set-TextValue "G22" "15"
Set-TextValue "D23" "0.1972"
Set-TextValue "E23" "10.18%"
Set-TextValue "G23" "15"
Set-TextValue "D24" "0.04534"
Set-TextValue "E24" "-0.90%"
Set-TextValue "G24" "15"
Set-TextValue "D25" "0.001230"
Set-TextValue "E25" "-1.46%"
Set-TextValue "G25" "15"
Set-TextValue "D26" "0.004163"
Set-TextValue "E26" "-7.79%"
Set-TextValue "G26" "15"
Set-TextValue "D27" "0.0001248"
Set-TextValue "E27" "0.09%"
Set-TextValue "G27" "15"
Set-TextValue "G28" "15"
Set-TextValue "G29" "15"
Set-TextValue "G30" "15"
Set-TextValue "G31" "15"
Set-TextValue "G32" "15"
Set-TextValue "G33" "15"
Set-TextValue "G34" "15"
Set-TextValue "G35" "15"
Set-TextValue "G36" "15"
Set-TextValue "G37" "15"
Set-TextValue "G38" "15"
Set-TextValue "D39" "0.01733"
Set-TextValue "E39" "-3.55%"
Set-TextValue "G39" "15"
Set-TextValue "D40" "0.04713"
Set-TextValue "E40" "-1.03%"
Set-TextValue "G40" "15"
Set-TextValue "D41" "0.007513"
Set-TextValue "E41" "4.03%"
Set-TextValue "G41" "15"
Set-TextValue "D42" "0.1364"
Set-TextValue "E42" "-0.83%"
Set-TextValue "G42" "15"
Set-TextValue "E43" "9.93%"
Set-TextValue "G43" "15"
Set-TextValue "D44" "0.01054"
Set-TextValue "E44" "6.77%"
Set-TextValue "G44" "15"
Set-TextValue "D45" "0.00006174"
Set-TextValue "E45" "-1.29%"
Set-TextValue "G45" "15"
Set-TextValue "E46" "0.07%"
Set-TextValue "G46" "15"
Set-TextValue "D47" "1.832"
Set-TextValue "E47" "62.95%"
Set-TextValue "G47" "15"
Set-TextValue "D48" "0.002995"
Set-TextValue "E48" "-16.63%"
Set-TextValue "G48" "15"
Set-TextValue "D49" "0.00002097"
Set-TextValue "E49" "0.07%"
Set-TextValue "G49" "15"
Set-TextValue "D50" "0.0001997"
Set-TextValue "E50" "0.07%"
Set-TextValue "G50" "15"
Set-TextValue "G51" "15"

Write-Host "Updated symbol list values (Price, Volume(1h), Hora) for rows 2-51"
